$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $c1 = $t.Cell($r, 1).Range.Text
    $c2 = $t.Cell($r, 2).Range.Text
    Write-Host ($r.ToString() + " | [" + $c1 + "] | [" + $c2 + "]")
}
